$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Proof-read / amended outcome (event summary) text for each of the four
# possible actions taken against the Military Commander during his trial.
$ws.Range("B15").Value = "It turns out that the Military Commander had been selling inside news to the neighbouring kingdom, providing them with information to form an attack plan on your kingdom."
$ws.Range("B16").Value = "The neighbouring kingdom had secretly infiltrated your kingdom in order to save the commander due to a pact he made with them. This angered you and you decide to wage war on the neighbouring kingdom."
$ws.Range("B17").Value = "The Military Commander was angered by this decision which made him outright betray the kingdom, joining the neighbouring kingdom. The neighbouring kingdom saw this as an opportunity to attack while you were without a Military Commander."
$ws.Range("B18").Value = "A sense of anger and betrayal was felt by the Military Commander due to this decision. His intention to sell critical kingdom secrets was intercepted by you. You made the decision to execute him for treason."
$ws.Range("B19").Value = "Putting the Military Commander on trial made him confess his deeds of selling information to the neighbouring kingdoms. This was deemed as treason and he was ultimately sent for execution."
